# Updates market-price derived columns (H..N) on several leve-profit rows
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets, as produced
# by the scheduled market-data refresh runner.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------- ALC ----
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H80").Value = 334.65
$ws.Range("J80").Value = 649.5
$ws.Range("L80").Value = 1948.5
$ws.Range("N80").Value = -3944.5

$ws.Range("H83").Value = 334.65
$ws.Range("J83").Value = 649.5
$ws.Range("L83").Value = 5845.5
$ws.Range("N83").Value = -15829.5

$ws.Range("H86").Value = 26049.35
$ws.Range("I86").Value = 33827
$ws.Range("J86").Value = 21861.385
$ws.Range("K86").Value = 33827
$ws.Range("L86").Value = 21861.385
$ws.Range("M86").Value = -32704
$ws.Range("N86").Value = -24107.385

$ws.Range("H88").Value = 26083098
$ws.Range("J88").Value = 30429780
$ws.Range("L88").Value = 30429780
$ws.Range("N88").Value = -30430592

$ws.Range("H89").Value = 26049.35
$ws.Range("I89").Value = 33827
$ws.Range("J89").Value = 21861.385
$ws.Range("K89").Value = 169135
$ws.Range("L89").Value = 109306.925
$ws.Range("M89").Value = -163519
$ws.Range("N89").Value = -120538.925

$ws.Range("H91").Value = 26083098
$ws.Range("J91").Value = 30429780
$ws.Range("L91").Value = 30429780
$ws.Range("N91").Value = -30432588

$ws.Range("H129").Value = 1039.1143
$ws.Range("I129").Value = 2250.3333
$ws.Range("J129").Value = 925.5625
$ws.Range("K129").Value = 6750.999899999999
$ws.Range("L129").Value = 2776.6875
$ws.Range("M129").Value = -1750.999899999999
$ws.Range("N129").Value = -12776.6875

# ---------------------------------------------------------------- ARM ----
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H32").Value = 1980.44
$ws.Range("I32").Value = 1980.0408
$ws.Range("K32").Value = 1980.0408
$ws.Range("M32").Value = -1693.0408

$ws.Range("H37").Value = 9878.857
$ws.Range("I37").Value = 8000
$ws.Range("J37").Value = 10192
$ws.Range("K37").Value = 8000
$ws.Range("L37").Value = 10192
$ws.Range("M37").Value = -7727
$ws.Range("N37").Value = -10738

$ws.Range("H43").Value = 5091.75
$ws.Range("J43").Value = 5091.75
$ws.Range("L43").Value = 5091.75
$ws.Range("N43").Value = -5717.75

$ws.Range("H55").Value = 10000
$ws.Range("I55").Value = 10000
$ws.Range("K55").Value = 10000
$ws.Range("M55").Value = -9685

$ws.Range("H63").Value = 1999.6666
$ws.Range("I63").Value = 1999.6666
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 1999.6666
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = ""
$ws.Range("N63").Value = -1313.6666

$ws.Range("H66").Value = 1999.6666
$ws.Range("I66").Value = 1999.6666
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 9998.333000000001
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = ""
$ws.Range("N66").Value = -6566.333000000001

$ws.Range("H132").Value = 18577386
$ws.Range("I132").Value = 20409366
$ws.Range("J132").Value = 7356503.5
$ws.Range("K132").Value = 61228098
$ws.Range("L132").Value = 22069510.5
$ws.Range("M132").Value = -61225568
$ws.Range("N132").Value = -22074570.5

# ---------------------------------------------------------------- BSM ----
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = ""
$ws.Range("N35").Value = 0

$ws.Range("H82").Value = 27283
$ws.Range("J82").Value = 27283
$ws.Range("L82").Value = 27283
$ws.Range("N82").Value = -28049

$ws.Range("H85").Value = 27283
$ws.Range("J85").Value = 27283
$ws.Range("L85").Value = 27283
$ws.Range("N85").Value = -29935

# ---------------------------------------------------------------- CRP ----
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H31").Value = 1303.3617
$ws.Range("I31").Value = 979.37933
$ws.Range("J31").Value = 1825.3334
$ws.Range("K31").Value = 979.37933
$ws.Range("L31").Value = 1825.3334
$ws.Range("M31").Value = -684.37933
$ws.Range("N31").Value = -2415.3334

$ws.Range("H34").Value = 1303.3617
$ws.Range("I34").Value = 979.37933
$ws.Range("J34").Value = 1825.3334
$ws.Range("K34").Value = 979.37933
$ws.Range("L34").Value = 1825.3334
$ws.Range("M34").Value = -777.37933
$ws.Range("N34").Value = -2229.3334

# ---------------------------------------------------------------- CUL ----
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H41").Value = 24480
$ws.Range("I41").Value = 500
$ws.Range("J41").Value = 60450
$ws.Range("K41").Value = 1500
$ws.Range("L41").Value = 181350
$ws.Range("M41").Value = -1162
$ws.Range("N41").Value = -182026

$ws.Range("H42").Value = 1000
$ws.Range("J42").Value = 1000
$ws.Range("L42").Value = 3000
$ws.Range("N42").Value = -4068

$ws.Range("H132").Value = 22733204
$ws.Range("I132").Value = 834
$ws.Range("J132").Value = 33341644
$ws.Range("K132").Value = 7506
$ws.Range("L132").Value = 300074796
$ws.Range("M132").Value = -4976
$ws.Range("N132").Value = -300079856

# ---------------------------------------------------------------- GSM ----
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H132").Value = 10276.381
$ws.Range("I132").Value = 7069.625
$ws.Range("J132").Value = 20538
$ws.Range("K132").Value = 21208.875
$ws.Range("L132").Value = 61614
$ws.Range("M132").Value = -18678.875
$ws.Range("N132").Value = -66674

# ---------------------------------------------------------------- LTW ----
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H46").Value = 1751.7333
$ws.Range("I46").Value = 612.8570999999999
$ws.Range("J46").Value = 2748.25
$ws.Range("K46").Value = 612.8570999999999
$ws.Range("L46").Value = 2748.25
$ws.Range("M46").Value = -424.8570999999999
$ws.Range("N46").Value = -3124.25

$ws.Range("H55").Value = 4968.2856
$ws.Range("I55").Value = 9293.817999999999
$ws.Range("J55").Value = 210.2
$ws.Range("K55").Value = 9293.817999999999
$ws.Range("L55").Value = 210.2
$ws.Range("M55").Value = -9120.817999999999
$ws.Range("N55").Value = -556.2

$ws.Range("H136").Value = 43419336
$ws.Range("I136").Value = 14882887
$ws.Range("K136").Value = 44648661
$ws.Range("M136").Value = -44646111

# ---------------------------------------------------------------- WVR ----
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H132").Value = 20515.322
$ws.Range("I132").Value = 24221.023
$ws.Range("J132").Value = 9645.267
$ws.Range("K132").Value = 72663.069
$ws.Range("L132").Value = 28935.801
$ws.Range("M132").Value = -70133.069
$ws.Range("N132").Value = -33995.801

$ws.Range("H136").Value = 10209053
$ws.Range("I136").Value = 13519029
$ws.Range("J136").Value = 3295.8333
$ws.Range("K136").Value = 40557087
$ws.Range("L136").Value = 9887.499899999999
$ws.Range("M136").Value = -40554537
$ws.Range("N136").Value = -14987.4999
